# JobMaterial.xlsx — "Digital" re-run of the Material sheet.
# Flexo/UV process rows are replaced with Digital Print equivalents; a new
# Film row (AC505T thermal transfer laminate) is added, the ink rows pick up
# new Digital ink codes (and reorder Yellow/Cyan), and the BOPP roll gets an
# updated quantity + a fresh (8th) row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$value)
    $range.NumberFormat = "@"
    $range.Value = $value
}

# Row 2 — Film / AC505T thermal transfer laminate (new material swapped in)
$ws.Range("C2").Value = "Digital Print 4x0"
$ws.Range("D2").Value = "AC505T THERMAL TRANSFER LAMINATE - 0.19 pt 11.5`" 10410 ppi"
Set-TextValue $ws.Range("E2") "555.36"
$ws.Range("H2").Value = ""

# Row 3 — Film / MATWA self-wound (carried down from the old row 2)
$ws.Range("A3").Value = "Film"
$ws.Range("C3").Value = "Digital Print 4x0"
$ws.Range("D3").Value = "MATWA 1M MATTE SELF WOUND - 1.52 pt 11.25`" 1310 ppi"
Set-TextValue $ws.Range("E3") "543.29"
$ws.Range("F3").Value = "M inch²"
$ws.Range("H3").Value = "37461 - MATWA 1M MATTE SELF WOUND - 1.52 pt  1310 ppi"

# Row 4 — Ink / Varnish / Black - Digital
$ws.Range("C4").Value = "Digital Print 4x0"
$ws.Range("D4").Value = "Black - Digital - "
Set-TextValue $ws.Range("E4") "0.50"
$ws.Range("H4").Value = "000015835 - Saphria Ink Black Labelfire UV02`nHT.400.1025/`n10L Bag In Box`n2.2lbs/Litre`n`$830.00/container"

# Row 5 — Ink / Varnish / Yellow - Digital
$ws.Range("C5").Value = "Digital Print 4x0"
$ws.Range("D5").Value = "Yellow - Digital - "
Set-TextValue $ws.Range("E5") "0.50"
$ws.Range("H5").Value = "000015833 - Saphria Ink Yellow Labelfire UV02`nHT.400.1025/`n10L Bag In Box`n2.2lbs/Litre`n`$830.00/container"

# Row 6 — Ink / Varnish / Cyan - Digital
$ws.Range("C6").Value = "Digital Print 4x0"
$ws.Range("D6").Value = "Cyan - Digital - "
Set-TextValue $ws.Range("E6") "0.50"
$ws.Range("H6").Value = "000015837 - Saphria Ink Cyan Labelfire UV02`nHT.400.1025/`n10L Bag In Box`n2.2lbs/Litre`n`$830.00/container"

# Row 7 — Ink / Varnish / Magenta - Digital (old Roll/BOPP row repurposed)
$ws.Range("A7").Value = "Ink / Varnish"
$ws.Range("C7").Value = "Digital Print 4x0"
$ws.Range("D7").Value = "Magenta - Digital - "
Set-TextValue $ws.Range("E7") "0.50"
$ws.Range("F7").Value = "lbs"
$ws.Range("H7").Value = "000015831 - Saphria Ink Magenta Labelfire UV02`nHT.400.1025/`n10L Bag In Box`n2.2lbs/Litre`n`$830.00/container"

# Row 8 — new Roll / BOPP WHITE row (updated quantity), styled to match
# the rest of the sheet's red 14pt Calibri cell look.
$row8 = $ws.Range("A8:H8")
$row8.Font.Name = "Calibri"
$row8.Font.Size = 14
$row8.Font.ColorIndex = 1

$ws.Range("A8").Value = "Roll"
$ws.Range("B8").Value = "Label  2p"
$ws.Range("C8").Value = "Digital Print 4x0"
$ws.Range("D8").Value = "BOPP WHITE - 5.28 pt 11.875`" 380 ppi"
Set-TextValue $ws.Range("E8") "660.73"
$ws.Range("F8").Value = "M inch²"
$ws.Range("G8").Value = ""
$ws.Range("H8").Value = ""
